# Commit: "Added history starting from 2009-March-31 historical data to all stocks"
# This particular sheet (stock 0185 / HSSEB) was missing 9 trading days between
# 2019-11-15 and 2019-11-29. Insert 9 new rows at row 798 (shifting everything
# below down by 9) and fill them with the recovered daily OHLCV data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows above the current row 798 ("2019-11-29"), pushing the
# existing rows 798..870 down to 807..879.
$ws.Rows("798:806").Insert()

# New row data: row, timestamp(epoch seconds UTC), date, open, high, low, close, volume
$newRows = @(
    @(798, 1574035200, "2019-11-18", 0.84,  0.84,  0.825, 0.825, 145100),
    @(799, 1574121600, "2019-11-19", 0.83,  0.845, 0.82,  0.84,  544400),
    @(800, 1574208000, "2019-11-20", 0.83,  0.84,  0.82,  0.825, 335800),
    @(801, 1574294400, "2019-11-21", 0.82,  0.825, 0.795, 0.805, 692300),
    @(802, 1574380800, "2019-11-22", 0.8,   0.805, 0.785, 0.8,   594600),
    @(803, 1574640000, "2019-11-25", 0.8,   0.8,   0.775, 0.775, 732600),
    @(804, 1574726400, "2019-11-26", 0.785, 0.795, 0.78,  0.78,  225900),
    @(805, 1574812800, "2019-11-27", 0.78,  0.79,  0.775, 0.78,  1045600),
    @(806, 1574899200, "2019-11-28", 0.78,  0.795, 0.78,  0.795, 127200)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    # Columns B (date) and C (id "0185") look like numbers/dates to Excel's
    # type inference, so force literal text with a leading apostrophe -
    # exactly like typing '2019-11-18 / '0185 into the cell - otherwise the
    # date would collapse to a serial number and "0185" would lose its
    # leading zero.
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = "'" + $r[2]
    $ws.Cells.Item($rowNum, 3).Value = "'0185"
    $ws.Cells.Item($rowNum, 4).Value = "HSSEB"
    $ws.Cells.Item($rowNum, 5).Value = $r[3]
    $ws.Cells.Item($rowNum, 6).Value = $r[4]
    $ws.Cells.Item($rowNum, 7).Value = $r[5]
    $ws.Cells.Item($rowNum, 8).Value = $r[6]
    $ws.Cells.Item($rowNum, 9).Value = $r[7]
}
